$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update existing values, drop columns that no longer apply
$ws.Range("D2").Value2 = 4800
$ws.Range("E2").Value2 = 46
$ws.Range("F2").Value2 = 46
$ws.Range("G2").Value2 = 224
$ws.Range("H2").Value2 = 171
$ws.Range("I2").Value2 = 171
$ws.Range("K2").Value2 = 9156
$ws.Range("L2").Value2 = 957
$ws.Range("M2").Value2 = 8199
$ws.Range("N2").Value2 = 8199
$ws.Range("P2").Value2 = 504
$ws.Range("Q2").Value2 = 273
$ws.Range("R2").Value2 = 31
$ws.Range("S2").Value2 = -336
$ws.Range("T2").Value2 = 59
$ws.Range("U2").Value2 = 214
$ws.Range("V2").Value2 = 164
$ws.Range("W2").Value2 = 0.96
$ws.Range("X2").Value2 = 3.56
$ws.Range("Y2").Value2 = 2.08
$ws.Range("Z2").Value2 = 1.87
$ws.Range("AA2").Value2 = 11.67
$ws.Range("AB2").Value2 = 1555.01
$ws.Range("AC2").Value2 = 1698
$ws.Range("AD2").Value2 = 39.93
$ws.Range("AE2").Value2 = 83393
$ws.Range("AF2").Value2 = 0.8100000000000001
$ws.Range("AG2").Value2 = 1500
$ws.Range("AH2").Value2 = 2.21
$ws.Range("AI2").Value2 = 86.17
$ws.Range("AJ2").Value2 = 10080029
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3: update existing values, drop columns that no longer apply
$ws.Range("D3").Value2 = 5005
$ws.Range("E3").Value2 = 65
$ws.Range("F3").Value2 = 65
$ws.Range("G3").Value2 = 196
$ws.Range("H3").Value2 = 94
$ws.Range("I3").Value2 = 94
$ws.Range("K3").Value2 = 9293
$ws.Range("L3").Value2 = 1090
$ws.Range("M3").Value2 = 8202
$ws.Range("N3").Value2 = 8202
$ws.Range("P3").Value2 = 504
$ws.Range("Q3").Value2 = -54
$ws.Range("R3").Value2 = 196
$ws.Range("S3").Value2 = 14
$ws.Range("T3").Value2 = 242
$ws.Range("U3").Value2 = -297
$ws.Range("V3").Value2 = 326
$ws.Range("W3").Value2 = 1.31
$ws.Range("X3").Value2 = 1.87
$ws.Range("Y3").Value2 = 1.14
$ws.Range("Z3").Value2 = 1.01
$ws.Range("AA3").Value2 = 13.3
$ws.Range("AB3").Value2 = 1544.31
$ws.Range("AC3").Value2 = 928
$ws.Range("AD3").Value2 = 53.71
$ws.Range("AE3").Value2 = 83423
$ws.Range("AF3").Value2 = 0.6
$ws.Range("AG3").Value2 = 1200
$ws.Range("AH3").Value2 = 2.41
$ws.Range("AI3").Value2 = 126.12
$ws.Range("AJ3").Value2 = 10080029
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4: update existing values, drop columns that no longer apply
$ws.Range("D4").Value2 = 5073
$ws.Range("E4").Value2 = 163
$ws.Range("F4").Value2 = 163
$ws.Range("G4").Value2 = 400
$ws.Range("H4").Value2 = 318
$ws.Range("I4").Value2 = 318
$ws.Range("K4").Value2 = 9322
$ws.Range("L4").Value2 = 910
$ws.Range("M4").Value2 = 8412
$ws.Range("N4").Value2 = 8412
$ws.Range("P4").Value2 = 504
$ws.Range("Q4").Value2 = 758
$ws.Range("R4").Value2 = 237
$ws.Range("S4").Value2 = -297
$ws.Range("T4").Value2 = 255
$ws.Range("U4").Value2 = 503
$ws.Range("V4").Value2 = 138
$ws.Range("W4").Value2 = 3.21
$ws.Range("X4").Value2 = 6.27
$ws.Range("Y4").Value2 = 3.83
$ws.Range("Z4").Value2 = 3.42
$ws.Range("AA4").Value2 = 10.82
$ws.Range("AB4").Value2 = 1583.98
$ws.Range("AC4").Value2 = 3154
$ws.Range("AD4").Value2 = 17.41
$ws.Range("AE4").Value2 = 85553
$ws.Range("AF4").Value2 = 0.64
$ws.Range("AG4").Value2 = 1500
$ws.Range("AH4").Value2 = 2.73
$ws.Range("AI4").Value2 = 46.39
$ws.Range("AJ4").Value2 = 10080029
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5: update existing values, drop columns that no longer apply
$ws.Range("D5").Value2 = 5584
$ws.Range("E5").Value2 = 326
$ws.Range("F5").Value2 = 326
$ws.Range("G5").Value2 = 381
$ws.Range("H5").Value2 = 317
$ws.Range("I5").Value2 = 317
$ws.Range("K5").Value2 = 9405
$ws.Range("L5").Value2 = 847
$ws.Range("M5").Value2 = 8558
$ws.Range("N5").Value2 = 8558
$ws.Range("P5").Value2 = 504
$ws.Range("Q5").Value2 = 583
$ws.Range("R5").Value2 = -290
$ws.Range("S5").Value2 = -199
$ws.Range("T5").Value2 = 51
$ws.Range("U5").Value2 = 532
$ws.Range("V5").Value2 = 83
$ws.Range("W5").Value2 = 5.84
$ws.Range("X5").Value2 = 5.67
$ws.Range("Y5").Value2 = 3.73
$ws.Range("Z5").Value2 = 3.38
$ws.Range("AA5").Value2 = 9.9
$ws.Range("AB5").Value2 = 1617.57
$ws.Range("AC5").Value2 = 3143
$ws.Range("AD5").Value2 = 20.17
$ws.Range("AE5").Value2 = 87045
$ws.Range("AF5").Value2 = 0.73
$ws.Range("AG5").Value2 = 1650
$ws.Range("AH5").Value2 = 2.6
$ws.Range("AI5").Value2 = 51.21
$ws.Range("AJ5").Value2 = 10080029
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6: update existing values, drop columns that no longer apply
$ws.Range("D6").Value2 = 5579
$ws.Range("E6").Value2 = 249
$ws.Range("F6").Value2 = 249
$ws.Range("G6").Value2 = 436
$ws.Range("H6").Value2 = 350
$ws.Range("I6").Value2 = 350
$ws.Range("K6").Value2 = 9502
$ws.Range("L6").Value2 = 787
$ws.Range("M6").Value2 = 8716
$ws.Range("N6").Value2 = 8716
$ws.Range("P6").Value2 = 504
$ws.Range("Q6").Value2 = 495
$ws.Range("R6").Value2 = -208
$ws.Range("S6").Value2 = -172
$ws.Range("T6").Value2 = 56
$ws.Range("U6").Value2 = 439
$ws.Range("V6").Value2 = 77
$ws.Range("W6").Value2 = 4.46
$ws.Range("X6").Value2 = 6.28
$ws.Range("Y6").Value2 = 4.06
$ws.Range("Z6").Value2 = 3.71
$ws.Range("AA6").Value2 = 9.029999999999999
$ws.Range("AB6").Value2 = 1654.89
$ws.Range("AC6").Value2 = 3476
$ws.Range("AD6").Value2 = 13.21
$ws.Range("AE6").Value2 = 88643
$ws.Range("AF6").Value2 = 0.52
$ws.Range("AG6").Value2 = 1600
$ws.Range("AH6").Value2 = 3.49
$ws.Range("AI6").Value2 = 44.9
$ws.Range("AJ6").Value2 = 10080029

# Row 7: drop all period-end metrics, keep id/name columns
$ws.Range("D7:AI7").ClearContents()

# Row 8: drop all period-end metrics, keep id/name columns
$ws.Range("D8:AI8").ClearContents()

# Row 9: drop all period-end metrics, keep id/name columns
$ws.Range("D9:AI9").ClearContents()
